$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Rename "CBSSports" -> "CreateAccount" and populate it with 6 rows of
# sign-up/create-account test data (header row + 5 test cases).
$ws3.Name = "CreateAccount"

$ws3.Columns.Item(1).ColumnWidth = 19.498697916666668
$ws3.Columns.Item(2).ColumnWidth = 14.998697916666666
$ws3.Columns.Item(3).ColumnWidth = 18.998697916666668
$ws3.Columns.Item(4).ColumnWidth = 15.830729166666666

# Header row
$ws3.Range("A1").Value = "First Name"
$ws3.Range("B1").Value = "Last Name"
$ws3.Range("C1").Value = "Email"
$ws3.Range("D1").Value = "Password"

# Row 2 - missing First Name
$ws3.Range("A2").Value = "'"
$ws3.Range("B2").Value = "Mike"
$ws3.Range("C2").Value = "mike34@yahoo.com"
$ws3.Hyperlinks.Add($ws3.Range("C2"), "mailto:mike34@yahoo.com")
$ws3.Range("D2").Value = "asdf8970"

# Row 3 - missing Last Name
$ws3.Range("A3").Value = "Jonathon"
$ws3.Range("B3").Value = "'"
$ws3.Range("C3").Value = "mike34@yahoo.com"
$ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:mike34@yahoo.com")
$ws3.Range("D3").Value = "asdf8970"

# Row 4 - missing Email
$ws3.Range("A4").Value = "Jonathon"
$ws3.Range("B4").Value = "Mike"
$ws3.Range("C4").Value = "'"
$ws3.Range("D4").Value = "asdf8970"

# Row 5 - missing Password
$ws3.Range("A5").Value = "Jonathon"
$ws3.Range("B5").Value = "Mike"
$ws3.Range("C5").Value = "mike34@yahoo.com"
$ws3.Hyperlinks.Add($ws3.Range("C5"), "mailto:mike34@yahoo.com")
$ws3.Range("D5").Value = "'"

# Row 6 - all fields filled (happy path)
$ws3.Range("A6").Value = "Jonathon"
$ws3.Range("B6").Value = "Mike"
$ws3.Range("C6").Value = "mike34@yahoo.com"
$ws3.Hyperlinks.Add($ws3.Range("C6"), "mailto:mike34@yahoo.com")
$ws3.Range("D6").Value = "asdf8970"

# Make CreateAccount the active/visible sheet with A3 selected.
$ws3.Activate()
$ws3.Range("A3").Select()
